# Auto-generated Excel COM-interop script to update cryptos.xlsx values per diff
# Updates Price (D) and Volume(1h) (E) columns, and for two swapped coin pairs also Coin (B) and Link (C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.362.99"
$ws.Range("E2").Value = "  -2.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.85"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.38"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4656"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.09"
$ws.Range("E9").Value = "  -11.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07980"
$ws.Range("E10").Value = "  -5.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.016"
$ws.Range("E11").Value = "  -3.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.54"
$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.46"
$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.945"
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.149"
$ws.Range("E15").Value = "  -3.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.32"
$ws.Range("E17").Value = "  -3.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001035"
$ws.Range("E18").Value = "  -3.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06558"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.30"
$ws.Range("E20").Value = "  -3.84%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.493"
$ws.Range("E22").Value = "  -4.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.361.25"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.066.41"
$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.50"
$ws.Range("E27").Value = "  +2.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.74"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.070"
$ws.Range("E29").Value = "  -3.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.486"
$ws.Range("E30").Value = "  -4.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.16"
$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09479"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9522"
$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.448"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.584"
$ws.Range("E35").Value = "  -1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.275"
$ws.Range("E36").Value = "  -5.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06060"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02234"
$ws.Range("E38").Value = "  -2.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.214"
$ws.Range("E39").Value = "  -4.08%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.055"
$ws.Range("E40").Value = "  -8.87%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5949"
$ws.Range("E42").Value = "  -3.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1897"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.26"
$ws.Range("E44").Value = "  -7.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.269"
$ws.Range("E45").Value = "  -3.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5654"
$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("E47").Value = "  -5.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.436"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  -4.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06770"
$ws.Range("E50").Value = "  -1.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.48"
$ws.Range("E51").Value = "  -1.54%  "
